$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-NumValue($ws, $addr, $val) {
    $ws.Range($addr).Value = $val
}

function Set-TextValue($ws, $addr, $val) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

# --- Cell value updates (rows 8-36) ---
Set-NumValue $ws 'C8' 81
Set-NumValue $ws 'C9' 2
Set-TextValue $ws 'G9' '512.00'
Set-NumValue $ws 'C10' 1
Set-TextValue $ws 'G10' '472.00'
Set-NumValue $ws 'C11' 75
Set-TextValue $ws 'G11' '49650.00'
Set-TextValue $ws 'A12' 'P. point'
Set-NumValue $ws 'C12' 39
Set-TextValue $ws 'D12' '6'
Set-TextValue $ws 'E12' 'On board'
Set-NumValue $ws 'F12' 136
Set-TextValue $ws 'G12' '5304.00'
Set-NumValue $ws 'C13' 82
Set-TextValue $ws 'G13' '1886.00'
Set-NumValue $ws 'C14' 42
Set-TextValue $ws 'G14' '2100.00'
Set-NumValue $ws 'C15' 41
Set-TextValue $ws 'D15' '6.0'
Set-TextValue $ws 'E15' 'Providing & Fixing of  3/6 pin 16 amp flush type non modular socket  made out from Industrial grade Polycarbonate or fire resistant ABS material, brass terminal with Porcelain based back cover & captive screws including cutting hole in tile and making connection testing etc. as required.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
Set-NumValue $ws 'F15' 78
Set-TextValue $ws 'G15' '3198.00'
Set-NumValue $ws 'C16' 76
Set-TextValue $ws 'D16' '8.0'
Set-TextValue $ws 'E16' 'Providing & Fixing of ISI marked (IS:1258) batten/angle lamp  holder with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material, brass terminal & captive screwsincluding making connection testing etc. as required.  All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
Set-NumValue $ws 'F16' 30
Set-TextValue $ws 'G16' '2280.00'
Set-NumValue $ws 'C17' 5
Set-TextValue $ws 'D17' '9.0'
Set-TextValue $ws 'E17' 'Providing & Fixing of IS 11037:1984  marked  non modular socket size flush type 180 watt rotary minimum 5 step fan regulator with CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including making connection testing etc. as required.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
Set-NumValue $ws 'F17' 219
Set-TextValue $ws 'G17' '1095.00'
Set-NumValue $ws 'C18' 24
Set-TextValue $ws 'D18' '10.0'
Set-TextValue $ws 'E18' 'Providing and fixing of   power plug point with non modular accessories as per PWD specification for electrical Works with  Galvanized   box of 1.2 mm thick  with earth terminal with suitable size phenolic laminated sheet (IS : 2036 -  1995) cover including cost of 16 amp. Switch (IS :3854) and 3/6 pin 16 amp. socket outlet  making connection , testing , etc. as required. . For specification of  Wiring accessories refer Chapter  E - 07 related item &  For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
Set-NumValue $ws 'F18' 303
Set-TextValue $ws 'G18' '7272.00'
Set-NumValue $ws 'C19' 73
Set-NumValue $ws 'C20' 64
Set-TextValue $ws 'G20' '2560.00'
Set-NumValue $ws 'C21' 44
Set-TextValue $ws 'G21' '2464.00'
Set-NumValue $ws 'C22' 71
Set-NumValue $ws 'C23' 2
Set-TextValue $ws 'D23' '19'
Set-TextValue $ws 'E23' '2 x 2.5 sq. mm. + 1x1.5sqmm'
Set-NumValue $ws 'F23' 81
Set-TextValue $ws 'G23' '162.00'
Set-TextValue $ws 'A24' 'Mtr.'
Set-NumValue $ws 'C24' 23
Set-TextValue $ws 'D24' '20'
Set-TextValue $ws 'E24' '2 x 4.0 sq. mm. + 1 x 2.5 sq. mm.'
Set-NumValue $ws 'F24' 122
Set-TextValue $ws 'G24' '2806.00'
Set-TextValue $ws 'A25' 'Set'
Set-NumValue $ws 'C25' 37
Set-TextValue $ws 'D25' '13.0'
Set-TextValue $ws 'E25' 'Plate Earthing  as per IS:3043 with Hot dipped G.I. Earth plate of size 600mm x 600mm x 6.0mm by embodying  3 to 4 mtr. below the ground level with 20  mm dia. G.I. ''B'' class watering Pipe ,including all accessories like nut, bolts, reducer, nipple, wire meshed funnel, and Heavy duty weather proof poly-propylene earth pit chamber with lockable Jam free lid suitable for safe working load 5000 Kg or more of size Top Dia. 225 to 260 mm, Bottom Dia 300 to 350 mm. and Height  250 to 300 mm. and embodying the pipe  complete with alternate layers salt and coke/ charcoal, testing of earth resistance for value of 5 ohms or less  as required & must record by engineer in charge during site visit and ensure to enter in measurment book.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .   '
Set-NumValue $ws 'F25' 5733
Set-TextValue $ws 'G25' '212121.00'
Set-TextValue $ws 'A26' ''
Set-NumValue $ws 'C26' 68
Set-TextValue $ws 'D26' '15.0'
Set-TextValue $ws 'E26' 'Providing & Fixing of  BEE  Star rated copper wounded double ball bearing capacitor start, aluminium body & Metallic  blade ceiling  fan  Conforming to all the performance requirements laid down in IS 374:2019 including all amendments, as applicable ; & Carry BIS licensing (i.e. ISI marking) with down rod up to 80 cm with secondary support safety cable ( steel rope) , cotter pin with 3 x 1.5 sq.mm pvc insulated flexible copper conductor making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
Set-NumValue $ws 'F26' 0
Set-TextValue $ws 'G26' '0.00'
Set-TextValue $ws 'A27' ''
Set-NumValue $ws 'C27' 31
Set-TextValue $ws 'D27' '16.0'
Set-TextValue $ws 'E27' 'Providing & Fixing of IP20 SMD Mid Power LED batten type integrated light fixture made from Powder coated Extruded aluminium  housing with in built driver  , System lumen efficacy ≥ 110 lm/Watt output, internal surge protection of 2.5 KV with Short & Open circuit protection ,THD < 10% , P. F.≥0.95, CRI >80 , life time of minimum  50000 Burning Hours with , 70% of intial Lumen maintaned till life ends  , CCT 3000°K / 4000°K  / 5700°K /6000°K/6500°K (As per ANSI Bin) , Maximum power consumption should not more than the specified rating and Fixture shall be of  BIS standard and  trade mark certificate ( T.C.). Manufactures Word Mark/ Name Engraved/ Embossing/ Screen printing on housing. OEM must have its own in house NABL lab setup for all testing facilities for LED fixtures. (LM79 & LM80) certificate / Report from OEM shall be submitted.  All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure "A" attached with this BSR .'
Set-NumValue $ws 'F27' 0
Set-TextValue $ws 'G27' '0.00'
Set-TextValue $ws 'A28' 'Each'
Set-NumValue $ws 'C28' 66
Set-TextValue $ws 'D28' '27'
Set-TextValue $ws 'E28' '1170mm(+/-10%) LED batten with min. lumen output 2200 lm'
Set-NumValue $ws 'F28' 492
Set-TextValue $ws 'G28' '32472.00'
Set-TextValue $ws 'A29' ''
Set-NumValue $ws 'C29' 48
Set-TextValue $ws 'D29' '17.0'
Set-TextValue $ws 'E29' 'Providing & Fixing of 240/415 V AC MCB with positive isolation of 10 kA breaking capacity (B/ C/D tripping characteristic as per type of load and  site requirement) 4 KV impulse withstand voltage, ISI marked IS 8828(1996) / conforming to IEC 60898-1 2002, IEC 60947-2, low watt losses, trip free mechanisum , energy limiting of  class 3 as per IEC,  minimum phase termination capacity of 35sq.mm. , conductor line load reversibility , IP 20 contact protection and fitted in  existing distribution board/sheets, minimum electrical operation 20,000 upto 20 A rating and 10,000 upto 63 A, 5000 for 80 A & above rating  including making connections, testing etc. as required. OEM shall have submit  NABL / CPRI / ERDA accrediated   lab type test reports  & All as per pre approved by Engineer in charge. For additional technical parameters of product / work refer Annexure ''A'' attached with this BSR'
Set-NumValue $ws 'F29' 0
Set-TextValue $ws 'G29' '0.00'
Set-NumValue $ws 'C30' 93
Set-TextValue $ws 'D30' '29'
Set-TextValue $ws 'E30' 'Single pole MCB   (With B/C curve tripping Characteristics)'
Set-TextValue $ws 'A31' 'Each'
Set-TextValue $ws 'D31' '30'
Set-TextValue $ws 'E31' ' 6 A to 32 A rating'
Set-NumValue $ws 'F31' 187
Set-TextValue $ws 'G31' '1122.00'
Set-NumValue $ws 'C32' 93
Set-NumValue $ws 'C33' 77
Set-TextValue $ws 'G33' '69300.00'
Set-NumValue $ws 'C34' 60
Set-NumValue $ws 'C35' 68
Set-NumValue $ws 'C36' 34
Set-TextValue $ws 'D36' '36'
Set-TextValue $ws 'E36' 'Total'

# --- Insert a new row at 37, shifting old rows 37-40 down to 38-41 ---
$ws.Rows("37").Insert()

# --- Populate the newly inserted row 37 ("Add Tender Premium") ---
Set-TextValue $ws 'A37' '%'
Set-NumValue $ws 'B37' 0
Set-NumValue $ws 'C37' 22
Set-TextValue $ws 'D37' '37'
Set-TextValue $ws 'E37' 'Add Tender Premium '
Set-NumValue $ws 'F37' 0
Set-TextValue $ws 'G37' '0.00'
Set-NumValue $ws 'H37' 0

# --- Update shifted totals rows (old 38->39, 39->40, 40->41) ---
Set-TextValue $ws 'G39' '396776.00'
Set-TextValue $ws 'H39' '396776.00'
Set-TextValue $ws 'G41' '396776.00'
Set-TextValue $ws 'H41' '396776.00'
